$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style prep: grow the font used by (the now-unused) G4 placeholder cell
# (cellXfs index 11) to Arial 13pt / #333333 so it can be reused by the new
# "Chassis Box" row's MFG# cell (A7).
$ws.Range("G4").Font.Name = "Arial"
$ws.Range("G4").Font.Size = 13
$ws.Range("G4").Font.Color = 3355443

# Move that formatting (and the now-empty cell) from G4 up to A7, then clear
# whatever is left behind in G4 so it goes back to being a bare cell.
$ws.Range("G4").Cut($ws.Range("A7"))
$ws.Range("G4").ClearFormats()

# --- Move the "Enclosure Parts" section header up from row 9 to row 6.
$ws.Range("A6").Value = "Enclosure Parts"
$ws.Range("A9").ClearContents()

# --- New row 7: Chassis Box part entry.
# Values are written in this particular order so that new shared-string
# table entries land at the same indices the original authoring session
# produced (description, then supplier, then the row-8 spacer, then the
# link, then the MFG# last).
$ws.Range("B7").Value = "NA"
$ws.Range("D7").Value = 12.59
$ws.Range("C7").Value = "Chassis Box - 1590XX, Diecast, 5.72"" x 4.77"" x 1.55"""
$ws.Range("E7").Value = "Amplified Parts"
$ws.Range("F7").Value = "Amplified Parts"

# --- New row 8: blank-ish spacer row with a single-space note under MFG.
$ws.Range("E8").Value = " "

$ws.Range("H7").Value = "https://www.amplifiedparts.com/products/chassis-box-1590xx-diecast-572-x-477-x-155"
$ws.Range("A7").Value = "P-H1590XXCE"
$ws.Rows(7).RowHeight = 16.5

# --- Selection, as left by the editing session.
$ws.Range("B8").Select()
